# "reduce tasks and metrics" - trims the Complexity column (G) entirely,
# collapses the per-row Dataset column (E) down to a constant "?" marker,
# and simplifies the Dialogue row's Quality Metric from "BLEU or Human" to "BLEU".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Cell values
# ---------------------------------------------------------------------

# Column E ("Dataset") becomes a flat "?" placeholder for every data row.
$ws.Range("E2").Value = "?"
$ws.Range("E3").Value = "?"
$ws.Range("E4").Value = "?"
$ws.Range("E5").Value = "?"
$ws.Range("E6").Value = "?"
$ws.Range("E7").Value = "?"

# Dialogue's quality metric shrinks from "BLEU or Human" to plain "BLEU".
$ws.Range("F6").Value = "BLEU"

# ---------------------------------------------------------------------
# 2. Re-style the cells whose look changes
#    (E1 now matches the other header cells; E2:E7 and F2:F7 now match
#    the plain "filled" body style used by column D)
# ---------------------------------------------------------------------

$ws.Range("A1").Copy()
$ws.Range("E1").PasteSpecial(-4122)

$ws.Range("D2:D7").Copy()
$ws.Range("E2:E7").PasteSpecial(-4122)
$ws.Range("F2:F7").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. Conditional formatting: drop every reference to columns E and G,
#    and remove the High/Medium/Low (Complexity) rule block entirely.
# ---------------------------------------------------------------------

# Block 1 ("No"/"Possible"/"Yes" rules) applied to C1:C1048576 E1:E1048576 G1:G1048576
# -> shrink to just C1:C1048576
$fcBlock1 = $ws.Range("E1:E1048576").FormatConditions.Item(1)
$fcBlock1.ModifyAppliesToRange($ws.Range("C1:C1048576"))

# Block 2 ("Yes" rule) applied to C2:C7 E2:E7 G2:G7 -> shrink to just C2:C7
$fcBlock2 = $ws.Range("E2:E7").FormatConditions.Item(1)
$fcBlock2.ModifyAppliesToRange($ws.Range("C2:C7"))

# Block 3 (High/Medium/Low, the Complexity rules) -> delete outright
$fcComplexity = $ws.Range("G1:G1048576")
while ($fcComplexity.FormatConditions.Count -gt 0) {
  $fcComplexity.FormatConditions.Item(1).Delete()
}

# ---------------------------------------------------------------------
# 4. Drop the Complexity column (G) altogether
# ---------------------------------------------------------------------

$ws.Columns.Item(7).Delete()

# ---------------------------------------------------------------------
# 5. Tidy up the remaining column widths to their new best-fit sizes
# ---------------------------------------------------------------------

$ws.Columns.Item(3).ColumnWidth = 11.166666666666666
$ws.Columns.Item(5).ColumnWidth = 7.0
$ws.Columns.Item(6).ColumnWidth = 58.0

# ---------------------------------------------------------------------
# 6. Selection follows the last touched cell
# ---------------------------------------------------------------------

$ws.Range("F7").Select()
